$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Account Maintenance")

# --- Update existing values ---
$ws.Range("D2").Value = "no"

# Fill new TestCaseId column (A) first ...
$ws.Range("A4").Value = "ACMO-0001"
$ws.Range("A5").Value = "ACMO-0002"
$ws.Range("A6").Value = "ACMO-0003"
$ws.Range("A7").Value = "ACMO-0004"

# ... then the new TestCaseDesc column (B) ...
$ws.Range("B4").Value = "Test case 1"
$ws.Range("B5").Value = "Test case 2"
$ws.Range("B6").Value = "Test case 3"
$ws.Range("B7").Value = "Test case 4"

# ... then the ExecutionStatus column (D) updates
$ws.Range("D4").Value = "YES"
$ws.Range("D6").Value = "YES"
$ws.Range("D7").Value = "YES"

# --- Fix formatting for A4 (style should match A5/A2 style - bordered, no fill) ---
$ws.Range("A5").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# --- Apply formatting to new rows 6 & 7: style 1 everywhere except col D (style 4) ---
$style1Src = $ws.Range("A5")
$style4Src = $ws.Range("D5")

foreach ($r in 6,7) {
    foreach ($col in @("A","B","C","E","F","G","H","I")) {
        $style1Src.Copy()
        $ws.Range("$col$r").PasteSpecial(-4122)
    }
    $style4Src.Copy()
    $ws.Range("D$r").PasteSpecial(-4122)
}
$ws.Application.CutCopyMode = $false

# --- Update selection to D3 ---
$ws.Range("D3").Select()
